# Applies the "Add files via upload" revision to Testing.xlsx:
#  - Adds a new "Modified Since Last Test" column (G) to every results table
#    on all three worksheets.
#  - Adds four new "Key invalid" validation test rows on the Validation sheet.
#  - Applies a dd/mm/yy date format to one (currently empty) Date Tested cell.
#  - Makes "Cogger Main" the active sheet / tab, and updates the remembered
#    cell selection on each sheet.

$wb = $excel.ActiveWorkbook

$newHeader = "Modified Since Last Test"

# ---------------------------------------------------------------------
# Sheet: Cogger Main
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Cogger Main")

$wsMain.Range("F1").Copy($wsMain.Range("G1"))
$wsMain.Range("G1").Value = $newHeader

# Touch the bottom-right corner so the sheet's used range grows to match
# the workbook this edit produced (extra blank padding rows/columns).
$wsMain.Cells.Item(19, 7).NumberFormat = "General"

$wsMain.Range("G1").Select()

# ---------------------------------------------------------------------
# Sheet: Validation
# ---------------------------------------------------------------------
$wsVal = $wb.Worksheets.Item("Validation")

# New column header on the first results table.
$wsVal.Range("F1").Copy($wsVal.Range("G1"))
$wsVal.Range("G1").Value = $newHeader

# Special punctuation used in the test-case names below: left/right
# double quotation marks and an en dash (built from char codes so the
# literal bytes in this script stay plain ASCII).
$ldq = [char]0x201C   # "
$rdq = [char]0x201D   # "
$endash = [char]0x2013   # -

# Insert the two new "At front - Single Wheel" / "At end - Three Wheels"
# rows for the "-" key-invalid tests (after row 10, before the old row 11).
$wsVal.Range("A11:A12").EntireRow.Insert()
$wsVal.Range("A11").Value = "Key invalid " + $ldq + $endash + $ldq + " At front " + $endash + " Single Wheel"
$wsVal.Range("A12").Value = "Key invalid " + $ldq + "-" + $rdq + " At end " + $endash + " Three Wheels"

# Insert the two new "At front - Single Wheel" / "At end - Three Wheels"
# rows for the "+" key-invalid tests (after row 14, before the old row 15).
$wsVal.Range("A15:A16").EntireRow.Insert()
$wsVal.Range("A15").Value = "Key invalid " + $ldq + "+" + $ldq + " At front " + $endash + " Single Wheel"
$wsVal.Range("A16").Value = "Key invalid " + $ldq + "+" + $rdq + " At end " + $endash + " Three Wheels"

# An (empty) Date Tested cell on the "Key Invalid Char - single Wheel" row
# picks up a dd/mm/yy date format.
$wsVal.Range("E6").NumberFormat = "dd/mm/yy"

# New column header on the second results table (the "Message Validation"
# block), which has shifted down to row 21 because of the 4 inserted rows.
$wsVal.Range("F21").Copy($wsVal.Range("G21"))
$wsVal.Range("G21").Value = $newHeader

# Touch the bottom-right corner so the sheet's used range grows to match
# the workbook this edit produced (extra blank padding rows/columns).
$wsVal.Cells.Item(25, 7).NumberFormat = "General"

$wsVal.Range("B22").Select()

# ---------------------------------------------------------------------
# Sheet: Encrypt and Decrypt
# ---------------------------------------------------------------------
$wsEnc = $wb.Worksheets.Item("Encrypt and Decrypt")

foreach ($headerRow in 1, 5, 8, 10) {
    $src = $wsEnc.Range("F" + $headerRow)
    $dst = $wsEnc.Range("G" + $headerRow)
    $src.Copy($dst)
    $dst.Value = $newHeader
}

# Touch the bottom-right corner so the sheet's used range grows to match
# the workbook this edit produced (extra blank padding rows/columns).
$wsEnc.Cells.Item(24, 7).NumberFormat = "General"

$wsEnc.Range("G10").Select()

# ---------------------------------------------------------------------
# Make "Cogger Main" the active sheet (activeTab = 0 in workbook.xml).
# ---------------------------------------------------------------------
$wsMain.Activate()
